$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.319.34"
$ws.Range("E2").Value = "  +6.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.695.44"
$ws.Range("E3").Value = "  +19.66%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.60"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.69"
$ws.Range("E6").Value = "  +7.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.693.67"
$ws.Range("E7").Value = "  +19.70%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +4.21%  "
$ws.Range("E10").Value = "  +8.33%  "
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("E12").Value = "  +5.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.83"
$ws.Range("E13").Value = "  +11.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("E14").Value = "  +6.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.309.91"
$ws.Range("E15").Value = "  +19.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.692.76"
$ws.Range("E16").Value = "  +19.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.342.36"
$ws.Range("E17").Value = "  +7.07%  "
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  +7.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.93"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "515.92"
$ws.Range("E21").Value = "  +6.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.23"
$ws.Range("E22").Value = "  +19.41%  "
$ws.Range("E23").Value = "  +8.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.43"
$ws.Range("E24").Value = "  +4.95%  "
$ws.Range("E25").Value = "  +8.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.42"
$ws.Range("E26").Value = "  +6.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  +8.81%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.52"
$ws.Range("E29").Value = "  +12.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.85"
$ws.Range("E31").Value = "  +14.32%  "
$ws.Range("E32").Value = "  +7.38%  "
$ws.Range("E33").Value = "  +18.42%  "
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.14"
$ws.Range("E36").Value = "  +10.27%  "
$ws.Range("E37").Value = "  +8.71%  "
$ws.Range("E38").Value = "  +10.62%  "
$ws.Range("E39").Value = "  +10.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.98"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("E41").Value = "  +4.42%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.176.44"
$ws.Range("E42").Value = "  +14.40%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.86"
$ws.Range("E43").Value = "  -7.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.79"
$ws.Range("E44").Value = "  +6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "409.95"
$ws.Range("E45").Value = "  +11.70%  "
$ws.Range("E46").Value = "  +6.80%  "
$ws.Range("E47").Value = "  +6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.17"
$ws.Range("E48").Value = "  +15.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.03"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.45"
$ws.Range("E51").Value = "  +13.22%  "
